$d = $word.ActiveDocument

# --- Locate the "Golang, " run -------------------------------------------
$golang = $d.Content.Duplicate
$found = $golang.Find.Execute("Golang, ", $true, $false, $false, $false, `
                               $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Golang, ' run to remove."
}

$gStart = $golang.Start
$gEnd   = $golang.End
$runLen = $gEnd - $gStart

$para      = $golang.Paragraphs(1).Range
$paraStart = $para.Start
$paraEnd   = $para.End

# --- Identify the neighbouring runs so the post-delete "merge adjacent ---
# --- runs with identical formatting" normalisation can't fuse them.   ---
# Immediately preceding run ("... R, Assembly, ").
$beforeStart = $gStart - 14
if ($beforeStart -lt $paraStart) { $beforeStart = $paraStart }
$before = $d.Range($beforeStart, $gStart)

# Immediately following run ("HTML").
$afterEnd = $gEnd + 4
if ($afterEnd -gt $paraEnd) { $afterEnd = $paraEnd }
$after = $d.Range($gEnd, $afterEnd)

# The two runs beyond that ("/CSS" then ", React") which would otherwise
# also get silently coalesced once the deletion ripples through them.
$after2End = $afterEnd + 4
if ($after2End -gt $paraEnd) { $after2End = $paraEnd }
$after2 = $d.Range($afterEnd, $after2End)

$after3 = $d.Range($after2End, $paraEnd)

# --- Remember each neighbour's original size, then bump every one of ---
# --- them to a distinct value so none of them match after the delete ---
$beforeSize = $before.Font.Size
$afterSize  = $after.Font.Size
$before.Font.Size = $beforeSize + 11
$after.Font.Size  = $afterSize + 12

$after2Size = $null
if ($after2.Start -lt $after2.End) {
    $after2Size = $after2.Font.Size
    $after2.Font.Size = $after2Size + 13
}

$after3Size = $null
if ($after3.Start -lt $after3.End) {
    $after3Size = $after3.Font.Size
    $after3.Font.Size = $after3Size + 14
}

# --- Remove the "Golang, " run entirely -----------------------------------
$d.Range($gStart, $gEnd).Delete()

# --- Restore original sizes, right-to-left (everything from $gEnd onward -
# --- shifted left by $runLen characters after the delete).           ---
if ($after3Size -ne $null) {
    $s = $after2End - $runLen
    $e = $paraEnd - $runLen
    $d.Range($s, $e).Font.Size = $after3Size
}

if ($after2Size -ne $null) {
    $s = $afterEnd - $runLen
    $e = $after2End - $runLen
    $d.Range($s, $e).Font.Size = $after2Size
}

$d.Range($gStart, $afterEnd - $runLen).Font.Size = $afterSize
$d.Range($beforeStart, $gStart).Font.Size = $beforeSize
